$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Clear the old summary tables (rows 2-4, cols L:U) and the old footnote
#    row (row 21) that are being relocated further down the sheet.
# ---------------------------------------------------------------------------
$ws.Range("L2:U4").ClearContents()
$ws.Range("L21:T21").ClearContents()

# ---------------------------------------------------------------------------
# 2. Re-create the "Facebook Graph Average" / "SNData Graph Average" summary
#    tables at their new location (rows 25-27, cols A:J) - the LSBFS results
#    were added as extra trial rows above (rows 5-7 already hold the data),
#    so the summary block now lives below all the raw samples.
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "Facebook Graph Average (7734 nodes)"
$ws.Range("F25").Value = "SNData Graph Average (124613 nodes)"

$ws.Range("A26").Value = "Sequential"
$ws.Range("B26").Value = "MTBFS"
$ws.Range("C26").Value = "LSBFS"
$ws.Range("D26").Value = "PPBFS"

$ws.Range("F26").Value = "Sequential"
$ws.Range("G26").Value = "MTBFS"
$ws.Range("H26").Value = "LSBFS"
$ws.Range("I26").Value = "PPBFS"

$ws.Range("A27").Formula = "=AVERAGE(A3:A7)"
$ws.Range("B27").Formula = "=AVERAGE(B3:B7)"
$ws.Range("C27").Formula = "=AVERAGE(C3:C7)"
$ws.Range("D27").Formula = "=AVERAGE(D3:D7)"
$ws.Range("E27").Value = "Time in ms"

$ws.Range("F27").Formula = "=AVERAGE(F3:F7)"
$ws.Range("G27").Formula = "=AVERAGE(G3:G7)"
$ws.Range("H27").Formula = "=AVERAGE(H3:H7)"
$ws.Range("I27").Formula = "=AVERAGE(I3:I7)"
$ws.Range("J27").Value = "Time in ms"

# ---------------------------------------------------------------------------
# 3. Re-create the footnote row at its new location (row 44).
# ---------------------------------------------------------------------------
$ws.Range("A44").Value = "Average time in milliseconds to run BFS on a FaceBook graph with 7734 nodes"
$ws.Range("I44").Value = "Average time in milliseconds to run BFS on scale-free graph with 124613 nodes"

# ---------------------------------------------------------------------------
# 4. Point the two charts at the relocated source ranges.
# ---------------------------------------------------------------------------
$chart1 = $ws.ChartObjects().Item(1)
$chart1.Chart.SeriesCollection(1).Formula = "=SERIES(,Sheet1!`$A`$26:`$D`$26,Sheet1!`$A`$27:`$D`$27,1)"

$chart2 = $ws.ChartObjects().Item(2)
$chart2.Chart.SeriesCollection(1).Formula = "=SERIES(,Sheet1!`$F`$26:`$I`$26,Sheet1!`$F`$27:`$I`$27,1)"

# ---------------------------------------------------------------------------
# 5. Move the charts themselves down/left to sit under the relocated tables.
# ---------------------------------------------------------------------------
$chart1.Left = 3.75
$chart1.Top = 423
$chart1.Width = 415.9345703125
$chart1.Height = 216

$chart2.Left = 520.0595703125
$chart2.Top = 422.25
$chart2.Width = 433.0625
$chart2.Height = 216

# ---------------------------------------------------------------------------
# 6. Update the sheet view: drop the frozen "topLeftCell=E1" scroll position
#    and move the active selection down one row (H16 -> H17).
# ---------------------------------------------------------------------------
$null = $ws.Range("H17").Select()
